$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "TOTAL SCORE" row entirely; every row below it shifts up by one.
$ws.Rows("22:22").Delete()

# The hyperlink that lived on the "Contract terms & conditions" comment cell
# needs to follow that row up from C26 to C25. The engine doesn't auto-shift
# the hyperlink's anchor on a row delete, so remove the stale one (still
# registered against the old "C26" address) and re-create it on the new
# cell - preserving the cell's existing text (the long eServices contract
# clause) instead of letting it get replaced by the link's display text.
$url = "https://content.vic.gov.au/sites/default/files/2023-12/eServcies-contract-%28April-2021%29.pdf"
$ws.Range("C26").Hyperlinks.Delete()

$linkCell = $ws.Range("C25")
$linkText = $linkCell.Text
$ws.Hyperlinks.Add($linkCell, $url, "", "", $url)
$linkCell.Value = $linkText

# Restore the previously-redacted evaluation comments referencing Wannon Water.
$ws.Range("C10").Value = "Higher score if company has previous worked with Wannon Water."
$ws.Range("C5").Value = "Higher score if company has previous worked with Wannon Water - this is important for continuity of delivery of Digital Operating Model and Digital Strategy.`nLower score if company has not previous worked for Wannon Water."

# Match the author's final cursor position.
$ws.Range("I6").Select()
